$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("device_recali_fields")

# Insert the new row just *below* the current row 2 instead of above it.
# Excel's default Insert() picks up formatting from the row above the
# insertion point, and the row directly above row 2 is the bold header row
# - inserting at row 3 instead means the copied formatting comes from row 2
# (a plain field row), which is what the new record should look like. The
# existing row 2 values are then shifted into the freshly formatted row 3,
# and the new cfres_id field is written into row 2.
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 2).Value = $ws.Cells.Item(2, 2).Value()
$ws.Cells.Item(3, 3).Value = $ws.Cells.Item(2, 3).Value()
$ws.Cells.Item(3, 4).Value = $ws.Cells.Item(2, 4).Value()

$ws.Cells.Item(2, 2).Value = "cfres_id"
$ws.Cells.Item(2, 3).Value = "string"
$ws.Cells.Item(2, 4).Value = "cfRes internal recall identifier"

# Match the row height used by sibling rows (17pt) for both the new row
# and the row that now holds the shifted-down event_date_initiated data.
$ws.Rows.Item(2).RowHeight = 17
$ws.Rows.Item(3).RowHeight = 17

# Update the visible selection to D2, as recorded after the edit.
$ws.Range("D2").Select() | Out-Null
